$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row index, Data (date serial), Interval, Prediction
$data = @(
    @(2, 45329, 12, 3.285614728927612),
    @(3, 45329, 13, 3.232341051101685),
    @(4, 45329, 14, 3.066033363342285),
    @(5, 45329, 15, 3.297250270843506),
    @(6, 45329, 16, 0.8487949967384338),
    @(7, 45329, 17, 0.01557119004428387),
    @(8, 45329, 18, 0.0154873700812459),
    @(9, 45329, 19, 0.02413088455796242),
    @(10, 45329, 20, 0.02413088455796242),
    @(11, 45329, 21, 0.02413088455796242),
    @(12, 45329, 22, 0.02413088455796242),
    @(13, 45329, 23, 0.02413088455796242),
    @(14, 45330, 0, 0.02396511659026146),
    @(15, 45330, 1, 0.02396511659026146),
    @(16, 45330, 2, 0.02396511659026146),
    @(17, 45330, 3, 0.02396511659026146),
    @(18, 45330, 4, 0.02396511659026146),
    @(19, 45330, 5, 0.02396511659026146),
    @(20, 45330, 6, 0.02396511659026146),
    @(21, 45330, 7, 0.02602160349488258),
    @(22, 45330, 8, 1.213412165641785),
    @(23, 45330, 9, 3.663368225097656),
    @(24, 45330, 10, 3.932492256164551),
    @(25, 45330, 11, 3.767731666564941),
    @(26, 45330, 12, 3.595699787139893),
    @(27, 45330, 13, 3.736693143844604),
    @(28, 45330, 14, 3.899230718612671),
    @(29, 45330, 15, 3.401404619216919),
    @(30, 45330, 16, 1.756325006484985),
    @(31, 45330, 17, 0.02620400488376617),
    @(32, 45330, 18, 0.02413088455796242),
    @(33, 45330, 19, 0.02413088455796242),
    @(34, 45330, 20, 0.02413088455796242),
    @(35, 45330, 21, 0.02413088455796242),
    @(36, 45330, 22, 0.02413088455796242),
    @(37, 45330, 23, 0.02413088455796242),
    @(38, 45331, 0, 0.02396511659026146),
    @(39, 45331, 1, 0.02396511659026146),
    @(40, 45331, 2, 0.02396511659026146),
    @(41, 45331, 3, 0.02396511659026146),
    @(42, 45331, 4, 0.02396511659026146),
    @(43, 45331, 5, 0.02396511659026146),
    @(44, 45331, 6, 0.02396511659026146),
    @(45, 45331, 7, 0.02602160349488258),
    @(46, 45331, 8, 1.450816631317139),
    @(47, 45331, 9, 3.631202459335327),
    @(48, 45331, 10, 3.944699287414551),
    @(49, 45331, 11, 3.767731666564941),
    @(50, 45331, 12, 3.595699787139893),
    @(51, 45331, 13, 3.71630072593689),
    @(52, 45331, 14, 3.913488149642944),
    @(53, 45331, 15, 3.45599889755249),
    @(54, 45331, 16, 1.756325006484985),
    @(55, 45331, 17, 0.02620400488376617),
    @(56, 45331, 18, 0.02413088455796242),
    @(57, 45331, 19, 0.02413088455796242),
    @(58, 45331, 20, 0.02413088455796242),
    @(59, 45331, 21, 0.02413088455796242),
    @(60, 45331, 22, 0.02413088455796242),
    @(61, 45331, 23, 0.02413088455796242),
    @(62, 45332, 0, 0.02396511659026146),
    @(63, 45332, 1, 0.02396511659026146),
    @(64, 45332, 2, 0.02396511659026146),
    @(65, 45332, 3, 0.02396511659026146),
    @(66, 45332, 4, 0.02396511659026146),
    @(67, 45332, 5, 0.02396511659026146),
    @(68, 45332, 6, 0.02396511659026146),
    @(69, 45332, 7, 0.02602160349488258),
    @(70, 45332, 8, 1.423936367034912),
    @(71, 45332, 9, 3.598320245742798),
    @(72, 45332, 10, 3.960421085357666),
    @(73, 45332, 11, 3.701785087585449),
    @(74, 45332, 12, 3.352742195129395),
    @(75, 45332, 13, 3.675429821014404),
    @(76, 45332, 14, 3.913488149642944),
    @(77, 45332, 15, 3.558686494827271),
    @(78, 45332, 16, 1.761136770248413),
    @(79, 45332, 17, 0.02620400488376617),
    @(80, 45332, 18, 0.02413088455796242),
    @(81, 45332, 19, 0.02413088455796242),
    @(82, 45332, 20, 0.02413088455796242),
    @(83, 45332, 21, 0.02413088455796242),
    @(84, 45332, 22, 0.02413088455796242),
    @(85, 45332, 23, 0.02413088455796242),
    @(86, 45333, 0, 3.173719882965088),
    @(87, 45333, 1, 3.173719882965088),
    @(88, 45333, 2, 3.234867572784424),
    @(89, 45333, 3, 3.234867572784424),
    @(90, 45333, 4, 3.173719882965088),
    @(91, 45333, 5, 3.268086194992065),
    @(92, 45333, 6, 3.234867572784424),
    @(93, 45333, 7, 3.153742074966431),
    @(94, 45333, 8, 3.207395076751709),
    @(95, 45333, 9, 3.410599708557129),
    @(96, 45333, 10, 3.423490047454834),
    @(97, 45333, 11, 3.223410844802856)
)

foreach ($item in $data) {
    $r = $item[0]
    $dateSerial = $item[1]
    $interval = $item[2]
    $prediction = $item[3]

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $interval
    $ws.Cells.Item($r, 3).Value = $prediction
}

Write-Host "Updated $($data.Count) rows"
